$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value. All target cells hold text-formatted
# numbers/strings (inlineStr in the source sheet), so every value is written as
# a quoted string and the cell is forced back to Text format beforehand so Excel
# does not silently reinterpret numeric-looking strings (e.g. "0.1380") as numbers
# and strip the significant trailing zeros.
$updates = [ordered]@{
    'D2' = '243.43'
    'D4' = '5.288'
    'D5' = '0.05783'
    'D7' = '3.334'
    'D8' = '0.8087'
    'D9' = '0.8788'
    'D10' = '0.1380'
    'D11' = '0.07300'
    'D12' = '0.03089'
    'D13' = '0.03059'
    'D14' = '0.09323'
    'D15' = '3.845'
    'D16' = '0.001543'
    'D17' = '0.04705'
    'D18' = '0.0006049'
    'D19' = '0.006115'
    'D20' = '0.001288'
    'D21' = '0.004601'
    'D22' = '0.00008695'
    'E22' = '21NitroExNTXBestin24h'
    'D23' = '3.581'
    'D28' = '0.0002350'
    'D40' = '0.03770'
    'D41' = '0.006426'
    'E41' = '40KickTokenKICK'
    'B42' = 'BKEXToken'
    'C42' = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
    'D42' = '0.1053'
    'E42' = '41BKEXTokenBKK'
    'B43' = 'CEJI'
    'C43' = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
    'D43' = '0.002419'
    'E43' = '42CEJICEJIWorstin24h'
    'D44' = '0.007132'
    'D45' = '0.00005472'
    'D47' = '0.5500'
    'D48' = '0.001856'
    'E48' = '47BOLOBOLO'
    'D49' = '0.00002100'
    'D50' = '0.0002000'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
